$d = $word.ActiveDocument

# Locate the (last) empty list paragraph that still only holds the
# "_GoBack" bookmark -- that's the paragraph the two new runs belong in.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "") {
        $target = $p
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Last
}

$newRunsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Name start with numbers, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>same as keywords</w:t></w:r></w:p>
'@

$insertAt = $target.Range.Start
$ins = $d.Range($insertAt, $insertAt)
$ins.InsertXML($newRunsXml)
